# Revert "Atualizacao dos .pdf e .pdf" on slide 1 (Telas-DIONISIO.pptx):
#   - shrink/move the subtitle placeholder back to its pre-update box
#   - drop the "10387644 - Fernando Karchiloff Gouveia de Amorim" line
#     that the reverted commit had appended

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle placeholder (ppPlaceholderSubtitle = 4) robustly,
# rather than hard-coding a shape index.
$subtitle = $null
$placeholders = $s.Shapes.Placeholders
for ($i = 1; $i -le $placeholders.Count; $i++) {
    $candidate = $placeholders.Item($i)
    if ($candidate.PlaceholderFormat.Type -eq 4) {
        $subtitle = $candidate
    }
}
if ($subtitle -eq $null) {
    $subtitle = $s.Shapes.Item(2)
}

# EMU -> point conversion (COM Left/Top/Width/Height are in points).
$emuPerPt = 12700

$subtitle.Left   = 1524000 / $emuPerPt
$subtitle.Top    = 3177687 / $emuPerPt
$subtitle.Width  = 9144000 / $emuPerPt
$subtitle.Height = 3255900 / $emuPerPt

# Remove the paragraph that reads
# "10387644 - Fernando Karchiloff Gouveia de Amorim" (the 7th paragraph,
# right after "10259046 - Samuel Silva Caetite").
$textRange = $subtitle.TextFrame.TextRange
$targetText = "10387644 - Fernando Karchiloff Gouveia de Amorim"

$paraCount = 1
$found = $false
while (-not $found) {
    $para = $textRange.Paragraphs($paraCount, 1)
    if ($para.Text.Length -eq 0) {
        break
    }
    if ($para.Text.TrimEnd("`r") -eq $targetText) {
        $para.Delete()
        $found = $true
    } else {
        $paraCount = $paraCount + 1
    }
}
